# EPBDS-9540 Support Java Name convention on Json field name generating in
# SpreadsheetResults.
#
# The "Uniquie cell name validation" test block (rows 4-6 on the active
# sheet) documents how raw step/result names are converted into Java/JSON
# field names. Update the sample names to reflect the new conversion rules
# and record the produced value ("1") next to the renamed steps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: expected field-name results for "Res_ult" and "Res".
$ws.Range("C4").Value = "ResUlt"
$ws.Range("D4").Value = "res"

# "ult_Value" step now converts to "ultValue"; its produced value is 1.
$ws.Range("B5").Value = "ultValue"
$ws.Range("C5").Value = "'1"

# "Value" step is unchanged, but now also records its produced value, 1.
$ws.Range("C6").Value = "'1"
